$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Helper: write a value to a cell, forcing TEXT storage when the
# string looks like a number (so it round-trips as inlineStr/string
# instead of being auto-coerced to a numeric cell).
# -----------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
    $ws.Cells.Item($row, $col).Value = $text
}

# ===================================================================
# Sheet "Alunos"
# ===================================================================
$ws1 = $wb.Worksheets.Item("Alunos")

# Row 2 updates
$ws1.Cells.Item(2, 6).Value = "['1', 't6959989', 't8890634']"
$ws1.Cells.Item(2, 7).Value = "['g2422223', 'g0330150']"
$ws1.Cells.Item(2, 8).Value = "[7.14, 5.69]"

# Row 3 updates
$ws1.Cells.Item(3, 6).Value = "['1', 't8890634']"
$ws1.Cells.Item(3, 7).Value = "['g0330150']"

# New row 6
Set-TextCell $ws1 6 1 "ra9392101"
Set-TextCell $ws1 6 2 "ra9392101"
$ws1.Cells.Item(6, 3).Value = "juao"
Set-TextCell $ws1 6 4 "23"
$ws1.Cells.Item(6, 5).Value = "aa@aa.com"
$ws1.Cells.Item(6, 6).Value = "[]"
$ws1.Cells.Item(6, 7).Value = "[]"

# New row 7
Set-TextCell $ws1 7 1 "ra8095505"
Set-TextCell $ws1 7 2 "ra8095505"
$ws1.Cells.Item(7, 3).Value = "daniel"
Set-TextCell $ws1 7 4 "20"
$ws1.Cells.Item(7, 5).Value = "bb@bb.com"
$ws1.Cells.Item(7, 6).Value = "[]"
$ws1.Cells.Item(7, 7).Value = "[]"

# ===================================================================
# Sheet "Turmas"
# ===================================================================
$ws2 = $wb.Worksheets.Item("Turmas")

$ws2.Cells.Item(2, 5).Value = "[{'id': '1', 'nome': 'C1', 'data_de_inicio': '12', 'data_de_fim': '21', 'peso_da_nota': '3'}, {'id': '2', 'nome': 'C2', 'data_de_inicio': '23', 'data_de_fim': '32', 'peso_da_nota': '4'}, {'id': 'c9878959', 'nome': 'C3', 'data_de_inicio': '12/12/2024', 'data_de_fim': '13/12/2024', 'peso_da_nota': '6'}, {'id': 'c4088151', 'nome': 'aba', 'data_de_inicio': '02/10/2024', 'data_de_fim': '01/11/2024', 'peso_da_nota': '5'}, {'id': 'c4238120', 'nome': 'c4', 'data_de_inicio': '27/11/2023', 'data_de_fim': '28/11/2023', 'peso_da_nota': '6'}, {'id': 'c1773595', 'nome': 'aaaaa', 'data_de_inicio': '27/11/2023', 'data_de_fim': '27/12/2023', 'peso_da_nota': '7'}, {'id': 'c5717717', 'nome': '3', 'data_de_inicio': '12/12/2023', 'data_de_fim': '12/12/2023', 'peso_da_nota': '6'}, {'id': 'c6884323', 'nome': 'c9878959', 'data_de_inicio': '28/11/2023', 'data_de_fim': '29/11/2023', 'peso_da_nota': '6'}, {'id': 'c5390957', 'nome': 'aaaaaaaa', 'data_de_inicio': '29/11/2023', 'data_de_fim': '30/11/2023', 'peso_da_nota': '7'}]"

# New row 7
Set-TextCell $ws2 7 1 "t8890634"
Set-TextCell $ws2 7 2 "t8890634"
$ws2.Cells.Item(7, 3).Value = "os batutinhas"
$ws2.Cells.Item(7, 4).Value = "31/01/2024"

# ===================================================================
# Sheet "Ciclos": insert 3 new rows before row 8 (pushes old row 8
# -> row 11), then populate rows 8-10 with the new cycles.
# ===================================================================
$ws3 = $wb.Worksheets.Item("Ciclos")
$ws3.Range("A8:A10").EntireRow.Insert()

Set-TextCell $ws3 8 1 "c5717717"
Set-TextCell $ws3 8 2 "c5717717"
Set-TextCell $ws3 8 3 "3"
Set-TextCell $ws3 8 4 "12/12/2023"
Set-TextCell $ws3 8 5 "12/12/2023"
Set-TextCell $ws3 8 6 "6"

Set-TextCell $ws3 9 1 "c6884323"
Set-TextCell $ws3 9 2 "c6884323"
$ws3.Cells.Item(9, 3).Value = "c9878959"
Set-TextCell $ws3 9 4 "28/11/2023"
Set-TextCell $ws3 9 5 "29/11/2023"
Set-TextCell $ws3 9 6 "6"

Set-TextCell $ws3 10 1 "c5390957"
Set-TextCell $ws3 10 2 "c5390957"
$ws3.Cells.Item(10, 3).Value = "aaaaaaaa"
Set-TextCell $ws3 10 4 "29/11/2023"
Set-TextCell $ws3 10 5 "30/11/2023"
Set-TextCell $ws3 10 6 "7"

# ===================================================================
# Sheet "Grupos"
# ===================================================================
$ws4 = $wb.Worksheets.Item("Grupos")

Set-TextCell $ws4 3 1 "g2391477"
Set-TextCell $ws4 3 2 "g2391477"
$ws4.Cells.Item(3, 3).Value = "hobbit's house"
$ws4.Cells.Item(3, 4).Value = "['1', '1', '1', '1', '1', '1', '2']"

Set-TextCell $ws4 4 1 "3"
Set-TextCell $ws4 4 2 "3"
$ws4.Cells.Item(4, 3).Value = "daaaa"
$ws4.Cells.Item(4, 4).Value = "['1', '1', '1', '1', '1', '1', '1', '2']"

Set-TextCell $ws4 5 1 "g2422223"
Set-TextCell $ws4 5 2 "g2422223"
$ws4.Cells.Item(5, 3).Value = "aaa"
$ws4.Cells.Item(5, 4).Value = "['1']"

Set-TextCell $ws4 6 1 "g0330150"
Set-TextCell $ws4 6 2 "g0330150"
$ws4.Cells.Item(6, 3).Value = "faa"
$ws4.Cells.Item(6, 4).Value = "['1', '2']"

# ===================================================================
# Sheet "Notas"
# ===================================================================
$ws5 = $wb.Worksheets.Item("Notas")

$ws5.Cells.Item(2, 5).Value = 5

Set-TextCell $ws5 4 1 "ID3"
Set-TextCell $ws5 4 2 "1"
$ws5.Cells.Item(4, 3).Value = "c1773595"
Set-TextCell $ws5 4 4 "1"
$ws5.Cells.Item(4, 5).Value = 4
